# Slide 2 ("1차 데모 시나리오"), content placeholder shape (Shapes.Item(2)):
#   - Paragraph 5 ("웹 : ...앱 갱신 + 반응형 웹") loses its trailing
#     "+ 반응형 웹" (and the space after "앱 갱신").
#   - Paragraph 6 ("푸시 알림 : db 값을 조회해서...") is reworded to
#     "푸시 알림 : DB를 조회해서 앱이 설치된 핸드폰에 알림 ".
#   - A new paragraph 7 "+ 터치 시 반응형 웹에 링크" is appended, carrying
#     the "반응형 웹" idea that used to live at the end of paragraph 5.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 6: "푸시 알림 : db 값을 조회해서 앱이 설치된 핸드폰에 알림을 한번 전달"
#     -> "푸시 알림 : DB를 조회해서 앱이 설치된 핸드폰에 알림 "
$para6 = $tr.Paragraphs(6)
# Reword the closing run while its index (5) is still valid.
$para6.Runs(5).Text = "를 조회해서 앱이 설치된 핸드폰에 알림 "
# Fold "db" + " " into the ": " run so it reads ": DB".
$para6.Runs(2).Text = ": DB"
# Drop the now-redundant "db" / " " runs (delete the later index first).
$para6.Runs(4).Text = ""
$para6.Runs(3).Text = ""

# --- New paragraph 7, right after paragraph 6: "+ 터치 시 반응형 웹에 링크"
$tr.Paragraphs(6).InsertAfter("`r+ ")
$para7 = $tr.Paragraphs(7)
$para7.InsertAfter("터치 시 반응형 웹에 링크")

# --- Paragraph 5: "웹 : ...웹&앱 갱신 + 반응형 웹" -> "웹 : ...웹&앱 갱신"
$para5 = $tr.Paragraphs(5)
$para5.Runs(5).Text = "앱 갱신"
$para5.Runs(7).Text = ""
$para5.Runs(6).Text = ""
